$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 684.4
$ws.Range("I135").Value = 649.3333
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 5843.9997
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -3308.9997
$ws.Range("N135").Value = -14070
$ws.Range("H137").Value = 1835.375
$ws.Range("I137").Value = 1271.2
$ws.Range("J137").Value = 2399.55
$ws.Range("K137").Value = 3813.6
$ws.Range("L137").Value = 7198.650000000001
$ws.Range("M137").Value = -1263.6
$ws.Range("N137").Value = -12298.65
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2702.389
$ws.Range("I61").Value = 2657.2307
$ws.Range("J61").Value = 2819.8
$ws.Range("K61").Value = 2657.2307
$ws.Range("L61").Value = 2819.8
$ws.Range("M61").Value = -2445.2307
$ws.Range("N61").Value = -3243.8
$ws.Range("H74").Value = 1512.9333
$ws.Range("I74").Value = 1319.1818
$ws.Range("K74").Value = 1319.1818
$ws.Range("M74").Value = -445.1818000000001
$ws.Range("H77").Value = 1512.9333
$ws.Range("I77").Value = 1319.1818
$ws.Range("K77").Value = 6595.909000000001
$ws.Range("M77").Value = -2227.909000000001
$ws.Range("H136").Value = 2702.389
$ws.Range("I136").Value = 2657.2307
$ws.Range("J136").Value = 2819.8
$ws.Range("K136").Value = 7971.6921
$ws.Range("L136").Value = 8459.400000000001
$ws.Range("M136").Value = -5421.6921
$ws.Range("N136").Value = -13559.4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H134").Value = 2014.4722
$ws.Range("I134").Value = 1656.1072
$ws.Range("J134").Value = 3268.75
$ws.Range("K134").Value = 4968.321599999999
$ws.Range("L134").Value = 9806.25
$ws.Range("M134").Value = -2433.321599999999
$ws.Range("N134").Value = -14876.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2336.1667
$ws.Range("I31").Value = 1950.2
$ws.Range("J31").Value = 2876.52
$ws.Range("K31").Value = 1950.2
$ws.Range("L31").Value = 2876.52
$ws.Range("M31").Value = -1655.2
$ws.Range("N31").Value = -3466.52
$ws.Range("H34").Value = 2336.1667
$ws.Range("I34").Value = 1950.2
$ws.Range("J34").Value = 2876.52
$ws.Range("K34").Value = 1950.2
$ws.Range("L34").Value = 2876.52
$ws.Range("M34").Value = -1748.2
$ws.Range("N34").Value = -3280.52
$ws.Range("H58").Value = 1483016.5
$ws.Range("I58").Value = 1684990
$ws.Range("J58").Value = 1878
$ws.Range("K58").Value = 1684990
$ws.Range("L58").Value = 1878
$ws.Range("M58").Value = -1684787
$ws.Range("N58").Value = -2284
$ws.Range("H132").Value = 324719.25
$ws.Range("I132").Value = 452980.38
$ws.Range("K132").Value = 1358941.14
$ws.Range("M132").Value = -1356411.14
$ws.Range("H134").Value = 2782.3
$ws.Range("I134").Value = 2102.875
$ws.Range("J134").Value = 5500
$ws.Range("K134").Value = 6308.625
$ws.Range("L134").Value = 16500
$ws.Range("M134").Value = -3773.625
$ws.Range("N134").Value = -21570
$ws.Range("H136").Value = 1483016.5
$ws.Range("I136").Value = 1684990
$ws.Range("J136").Value = 1878
$ws.Range("K136").Value = 5054970
$ws.Range("L136").Value = 5634
$ws.Range("M136").Value = -5052420
$ws.Range("N136").Value = -10734
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 250889.28
$ws.Range("I68").Value = 435384.7
$ws.Range("K68").Value = 1306154.1
$ws.Range("M68").Value = -1305343.1
$ws.Range("H71").Value = 250889.28
$ws.Range("I71").Value = 435384.7
$ws.Range("K71").Value = 3918462.3
$ws.Range("M71").Value = -3914406.3
$ws.Range("H112").Value = 5316
$ws.Range("J112").Value = 5929.6553
$ws.Range("L112").Value = 17788.9659
$ws.Range("N112").Value = -20004.9659
$ws.Range("H115").Value = 3950
$ws.Range("I115").Value = 2000
$ws.Range("J115").Value = 4166.6665
$ws.Range("K115").Value = 6000
$ws.Range("L115").Value = 12499.9995
$ws.Range("M115").Value = -4825
$ws.Range("N115").Value = -14849.9995
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1968.25
$ws.Range("I122").Value = 1968.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5904.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3454.75
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 5004.6665
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 6007
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 18021
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = -22961
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 500163
$ws.Range("J69").Value = 500163
$ws.Range("L69").Value = 500163
$ws.Range("N69").Value = -501785
$ws.Range("H72").Value = 500163
$ws.Range("J72").Value = 500163
$ws.Range("L72").Value = 1500489
$ws.Range("N72").Value = -1508601
$ws.Range("H74").Value = 15197
$ws.Range("I74").Value = 15197
$ws.Range("K74").Value = 15197
$ws.Range("M74").Value = -14199
$ws.Range("H77").Value = 15197
$ws.Range("I77").Value = 15197
$ws.Range("K77").Value = 45591
$ws.Range("M77").Value = -40599
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H132").Value = 5339.0356
$ws.Range("I132").Value = 5626.9473
$ws.Range("J132").Value = 4731.222
$ws.Range("K132").Value = 16880.8419
$ws.Range("L132").Value = 14193.666
$ws.Range("M132").Value = -14350.8419
$ws.Range("N132").Value = -19253.666
$ws.Range("H136").Value = 2720.2
$ws.Range("I136").Value = 1324.25
$ws.Range("J136").Value = 3650.8333
$ws.Range("K136").Value = 3972.75
$ws.Range("L136").Value = 10952.4999
$ws.Range("M136").Value = -1422.75
$ws.Range("N136").Value = -16052.4999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 39711.8
$ws.Range("J131").Value = 39711.8
$ws.Range("L131").Value = 39711.8
$ws.Range("N131").Value = -49791.8
$ws.Range("H132").Value = 4493.638
$ws.Range("I132").Value = 4523
$ws.Range("J132").Value = 4408
$ws.Range("K132").Value = 13569
$ws.Range("L132").Value = 13224
$ws.Range("M132").Value = -11039
$ws.Range("N132").Value = -18284
$ws.Range("H136").Value = 3728.8462
$ws.Range("I136").Value = 4345
$ws.Range("J136").Value = 3200.7144
$ws.Range("K136").Value = 13035
$ws.Range("L136").Value = 9602.143199999999
$ws.Range("M136").Value = -10485
$ws.Range("N136").Value = -14702.1432
